$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing Growth (H) values for rows 191-196 ---
$ws.Range("H191").Value2 = 0.5
$ws.Range("H192").Value2 = 1
$ws.Range("H193").Value2 = 0.75
$ws.Range("H194").Value2 = 1.5
$ws.Range("H195").Value2 = 3.5
$ws.Range("H196").Value2 = 9.5

# --- 2. Append 7 new data rows (562-568) ---
# Columns: A Date | B Plant_Type | C Plant_Size | D Low | E High | F Temp_Diff (formula)
# G Rain | H Growth | I Pruned | J Quadrant | K Shade | L UV | M Humidity | N Dew_Point
# O Pressure | P Wind_Gust | Q Cloud_Cover | R Visibility | S AQI | T Pollen

$newRows = @(
    @{ Row=562; B="Nonflowering"; C="Large";  D=73; E=88; G=5.57; H=1.1000000000000001; I="No"; J=2; K="Neutral"; L=6; M=0.66; N=75; O=30.05; P=4; Q=0.59; R=9.9; S=52; T=22 },
    @{ Row=563; B="Tree";         C="Medium"; D=73; E=88; G=5.57; H=1.2;                I="No"; J=3; K="Neutral"; L=6; M=0.66; N=75; O=30.05; P=4; Q=0.59; R=9.9; S=52; T=22 },
    @{ Row=564; B="Tree";         C="Small";  D=73; E=88; G=5.57; H=1.2;                I="No"; J=3; K="Dark";    L=6; M=0.66; N=75; O=30.05; P=4; Q=0.59; R=9.9; S=52; T=22 },
    @{ Row=565; B="Tree";         C="Medium"; D=73; E=88; G=5.57; H=1.5;                I="No"; J=3; K="Dark";    L=6; M=0.66; N=75; O=30.05; P=4; Q=0.59; R=9.9; S=52; T=22 },
    @{ Row=566; B="Tree";         C="Medium"; D=73; E=88; G=5.57; H=1.6;                I="No"; J=3; K="Neutral"; L=6; M=0.66; N=75; O=30.05; P=4; Q=0.59; R=9.9; S=52; T=22 },
    @{ Row=567; B="Tree";         C="Large";  D=73; E=88; G=5.57; H=4.5;                I="No"; J=4; K="Dark";    L=6; M=0.66; N=75; O=30.05; P=4; Q=0.59; R=9.9; S=52; T=22 },
    @{ Row=568; B="Pruned";       C="Medium"; D=73; E=88; G=5.57; H=11;                 I="No"; J=1; K="Bright";  L=6; M=0.66; N=75; O=30.05; P=4; Q=0.59; R=9.9; S=52; T=22 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value2 = 45867
    $ws.Range("A$row").NumberFormat = "m/d/yy"

    $ws.Range("B$row").Value2 = $r.B
    $ws.Range("C$row").Value2 = $r.C
    $ws.Range("D$row").Value2 = $r.D
    $ws.Range("E$row").Value2 = $r.E
    $ws.Range("G$row").Value2 = $r.G
    $ws.Range("H$row").Value2 = $r.H
    $ws.Range("I$row").Value2 = $r.I
    $ws.Range("J$row").Value2 = $r.J
    $ws.Range("K$row").Value2 = $r.K
    $ws.Range("L$row").Value2 = $r.L
    $ws.Range("M$row").Value2 = $r.M
    $ws.Range("N$row").Value2 = $r.N
    $ws.Range("O$row").Value2 = $r.O
    $ws.Range("P$row").Value2 = $r.P
    $ws.Range("Q$row").Value2 = $r.Q
    $ws.Range("R$row").Value2 = $r.R
    $ws.Range("S$row").Value2 = $r.S
    $ws.Range("T$row").Value2 = $r.T
}

# Fill the Temp_Diff shared formula down through the new rows (extends F543:F561 -> F543:F568)
$ws.Range("F543:F568").Formula = "=ABS(D543-E543)"

# --- 3. Restore selection to match the post-edit view state ---
$ws.Range("Q562:Q568").Select()
